# Generate Report for Handback
# The 7a3b1b1b-...md file has now been handed back, so its status moves
# from "Ready for handoff" to "Handed back: in sync with en-US" on every
# sheet, and the "Latest Handback DateTime" for both files is refreshed.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = $handedBack
$ws.Range("C3").Value = $handedBack

# --- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("B3").Value = $handedBack
$ws.Range("G2").Value = "2016-03-07 02:53:22"
$ws.Range("G3").Value = "2016-03-07 02:53:22"

# --- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("B3").Value = $handedBack
$ws.Range("G2").Value = "2016-03-07 02:53:42"
$ws.Range("G3").Value = "2016-03-07 02:53:42"
